# Add a new "Questions?" closing slide to the end of the deck.
#
# Matches the target edit: a new slide (the 15th) is appended using the
# "Animated Closing Slide" layout (CustomLayout index 32 - the layout whose
# only placeholder is the body/idx=10 "call to action" placeholder), with
# its placeholder text set to "Questions?" and a Fade slide transition.

$p = $ppt.ActivePresentation

# CustomLayout #32 on the slide master is "Animated Closing Slide", which
# exposes a single body placeholder (idx 10) - the same layout used by the
# new slide in the target presentation.
$layout = $p.SlideMaster.CustomLayouts.Item(32)

# Append the new slide after the current last slide (there are 14 slides,
# so the new one becomes slide 15).
$newIndex = $p.Slides.Count + 1
$s = $p.Slides.AddSlide($newIndex, $layout)

# Fill in the placeholder text.
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Questions?"

# Give the new slide the same Fade transition used throughout the deck.
$s.SlideShowTransition.EntryEffect = 1793  # ppEffectFade
